# "Generate Report for Handoff"
#
# The localization-status report is being regenerated: the status that was
# "Handed back: in sync with en-US" is now "Ready for handoff", and the
# timestamps that record when each language's handoff XLIFF was (re)generated
# move forward to reflect the new run.
#
# Status text "Handed back: in sync with en-US" -> "Ready for handoff"
# appears on every sheet (Overview!E2:F2, zh-cn!C2, de-de!C2).
#
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps are
# bumped to the new handoff run time. zh-cn gets its own, slightly earlier,
# timestamp (it was generated first); Overview and de-de share the later one.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status column: ready for handoff ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff-generation timestamps ---
$wsOverview.Range("G2").Value = "2016-09-04 17:05:18"
$wsZhCn.Range("H2").Value     = "2016-09-04 17:05:13"
$wsDeDe.Range("H2").Value     = "2016-09-04 17:05:18"

# The status text shrank quite a bit ("Handed back: in sync with en-US" ->
# "Ready for handoff"), so re-fit the columns that hold it so the report
# doesn't leave a ton of dead whitespace behind.
$wsOverview.Columns.Item(5).EntireColumn.AutoFit()
$wsOverview.Columns.Item(6).EntireColumn.AutoFit()
$wsZhCn.Columns.Item(3).EntireColumn.AutoFit()
$wsDeDe.Columns.Item(3).EntireColumn.AutoFit()
